$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.896.42'
$ws.Range('E2').Value = '  +0.52%  '
$ws.Range('D3').Value = '2.469.82'
$ws.Range('E3').Value = '  +0.39%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '560.76'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.08%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '163.64'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.38%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  +1.66%  '
$ws.Range('E9').Value = '  +4.05%  '
$ws.Range('E10').Value = '  +0.75%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.334'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.62%  '
$ws.Range('E12').Value = '  +0.82%  '
$ws.Range('D13').Value = '68.734.07'
$ws.Range('E13').Value = '  +0.55%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000171'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.38%  '
$ws.Range('E15').Value = '  +0.64%  '
$ws.Range('E16').Value = '  -3.15%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '337.17'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E18').Value = '  -3.50%  '
$ws.Range('E19').Value = '  -0.01%  '
$ws.Range('E21').Value = '  +0.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '66.67'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.75%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.66'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.96%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.27'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.09%  '
$ws.Range('D25').Value = '0.0₃0825'
$ws.Range('E25').Value = '  -1.75%  '
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '430.33'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.15%  '
$ws.Range('E29').Value = '  -2.19%  '
$ws.Range('E30').Value = '  -2.55%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '159.58'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.94%  '
$ws.Range('E32').Value = '  +0.02%  '
$ws.Range('E33').Value = '  -0.04%  '
$ws.Range('E34').Value = '  -1.30%  '
$ws.Range('E35').Value = '  -0.48%  '
$ws.Range('E36').Value = '  -0.57%  '
$ws.Range('E37').Value = '  -2.42%  '
$ws.Range('E38').Value = '  -3.52%  '
$ws.Range('E39').Value = '  -1.14%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.07'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.34%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.38'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.15%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '130.43'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.39%  '
$ws.Range('E43').Value = '  +0.20%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.485'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.25%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.566'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.34%  '
$ws.Range('E46').Value = '  +0.81%  '
$ws.Range('E47').Value = '  +0.22%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.39'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.57%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.00'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -7.12%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '16.90'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.29%  '
$ws.Range('D51').Value = '0.0₆0207'
$ws.Range('E51').Value = '  +0.80%  '
